$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.400.92"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "3.759.47"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'592.47"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "'165.26"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("D7").Value = "3.758.16"
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("D11").Value = "'6.37"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  -3.69%  "
$ws.Range("D14").Value = "'35.70"
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "4.394.25"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "3.744.94"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "67.419.98"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "'17.61"
$ws.Range("E18").Value = "  -3.71%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'6.89"
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").Value = "'10.48"
$ws.Range("E21").Value = "  -4.71%  "
$ws.Range("D22").Value = "'454.57"
$ws.Range("E22").Value = "  -3.05%  "
$ws.Range("D23").Value = "'0.692"
$ws.Range("E23").Value = "  -1.54%  "
$ws.Range("D24").Value = "'0.0000150"
$ws.Range("E24").Value = "  +5.38%  "
$ws.Range("D25").Value = "'83.22"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("E26").Value = "  -4.71%  "
$ws.Range("D27").Value = "'11.79"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "'2.77"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").Value = "'29.57"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").Value = "'7.16"
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("D33").Value = "'2.14"
$ws.Range("E33").Value = "  -3.27%  "
$ws.Range("D34").Value = "'9.12"
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "3.714.44"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "'0.0995"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("D40").Value = "'0.991"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "'5.71"
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'43.84"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").Value = "'0.298"
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").Value = "'46.76"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("D47").Value = "'1.88"
$ws.Range("E47").Value = "  -4.34%  "
$ws.Range("D48").Value = "'8.32"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").Value = "'146.37"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "'388.75"
$ws.Range("E50").Value = "  -4.41%  "
$ws.Range("D51").Value = "2.746.33"
$ws.Range("E51").Value = "  +2.24%  "
